$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# ---------------------------------------------------------------------------
# 1) A brand-new September transaction ("adani icici" @ 2024-09-15 13:10:50)
#    was logged. It lands at the top of the September_Details/September_Date
#    list (columns R/S, row 39) and pushes the existing 102 entries
#    (rows 39..140) down by one row, into rows 40..141 (row 141 was blank
#    before and now receives what used to be the last entry, row 140).
# ---------------------------------------------------------------------------

$firstRow = 39
$lastRow  = 140

# Snapshot the existing R/S values (top to bottom) before overwriting anything.
$snapshot = @()
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $rText = $ws.Cells.Item($row, 18).Text
    $sText = $ws.Cells.Item($row, 19).Text
    $snapshot += ,@($rText, $sText)
}

# Write the snapshot back shifted down by one row (bottom-up so we never
# clobber a source row before it has been read -- already safe since we
# fully captured the snapshot above, but keep bottom-up for clarity/safety).
for ($i = $snapshot.Count - 1; $i -ge 0; $i--) {
    $destRow = $firstRow + 1 + $i
    $ws.Cells.Item($destRow, 18).Value = $snapshot[$i][0]
    $ws.Cells.Item($destRow, 19).Value = $snapshot[$i][1]
}

# New top entry.
$ws.Cells.Item($firstRow, 18).Value = "adani icici"
$ws.Cells.Item($firstRow, 19).Value = "2024-09-15 13:10:50"

# ---------------------------------------------------------------------------
# 2) The "Broadband" group label moves from A149 down to the newly appended
#    A150, growing the sheet's used range to A1:Y150.
# ---------------------------------------------------------------------------

$ws.Range("A149").Value = ""
$ws.Range("A150").Value = "Broadband"
